$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Spring 1")
$ws.Activate()

# "Avance" (progress) column for the Spring 1 burndown table/chart.
$ws.Range("C2").Value = 5
$ws.Range("C3").Value = 4
$ws.Range("C4").Value = 4
$ws.Range("C5").Value = 3
$ws.Range("C6").Value = 2
$ws.Range("C7").Value = 0

$ws.Range("D9").Select()
